$wb = $excel.ActiveWorkbook

# Sheet 1 (s__CAG-314 sp000437915-b-p): remove the two rows for
# label_UMGS475_0.fasta and label_UMGS475_12.fasta (originally rows 5-6),
# shifting the remaining rows up. Dimension A1:E19 -> A1:E17.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A5:A6").EntireRow.Delete()

# Sheet 2 (s__CAG-314 sp900551395-b-p): remove the 20-row block
# (originally rows 21-40), shifting the remaining rows up.
# Dimension A1:E117 -> A1:E97.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A21:A40").EntireRow.Delete()
